$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "yes"
$ws.Range("C2").Value = "na"
$ws.Range("D2").Value = "na"

$ws.Range("B3").Select()
